$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MT"
$ws.Range("C2").Value = "Exp"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.9922
$ws.Range("F2").Value = 0.74
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2637.9293
$ws.Range("I2").Value = 0.2217
